# Regenerate save_data "K" column (strikeouts) values to reflect the
# recalculated figures from the (re-pulled) box-score data.
# Only column G (header "K") changes; all other columns/rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value
$updates = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 2
    8  = 1
    10 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
